# Regenerate merged AHB files
# 1. Rename header cells: "_old" -> "_FV2404" and "_new" -> "_FV2410"
# 2. Convert the used range into an Excel Table (ListObject) named "Table1"
# 3. Freeze the header row (pane split) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row cells --------------------------------------
$lastCol = 21  # columns A..U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = [string]$cell.Value2
    if ($text -like "*_old") {
        $cell.Value = $text -replace "_old$", "_FV2404"
    } elseif ($text -like "*_new") {
        $cell.Value = $text -replace "_new$", "_FV2410"
    }
}

# --- 2. Create the Excel Table over the used range ------------------------
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range("A1:U" + $lastRow)

$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
